$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "64.206.33"
Set-TextValue "E2" "  +0.18%  "
Set-TextValue "D3" "3.120.31"
Set-TextValue "E3" "  -8.05%  "
Set-TextValue "E4" "  -1.68%  "
Set-TextValue "D5" "590.80"
Set-TextValue "E5" "  +0.53%  "
Set-TextValue "D6" "153.84"
Set-TextValue "E6" "  +5.15%  "
Set-TextValue "E7" "  -1.43%  "
Set-TextValue "D8" "3.120.21"
Set-TextValue "E8" "  -0.31%  "
Set-TextValue "E9" "  +0.66%  "
Set-TextValue "E10" "  -0.01%  "
Set-TextValue "D11" "5.97"
Set-TextValue "E11" "  +2.15%  "
Set-TextValue "E12" "  +0.96%  "
Set-TextValue "B13" "ShibaInu"
Set-TextValue "C13" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D13" "0.0000246"
Set-TextValue "E13" "  -0.38%  "
Set-TextValue "B14" "Avalanche"
Set-TextValue "C14" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D14" "38.03"
Set-TextValue "E14" "  +2.37%  "
Set-TextValue "D15" "3.633.64"
Set-TextValue "E15" "  -7.60%  "
Set-TextValue "E16" "  -1.52%  "
Set-TextValue "D17" "7.23"
Set-TextValue "D18" "64.024.24"
Set-TextValue "E18" "  +0.15%  "
Set-TextValue "D19" "3.121.01"
Set-TextValue "E19" "  -2.01%  "
Set-TextValue "D20" "473.84"
Set-TextValue "E20" "  +2.18%  "
Set-TextValue "D21" "14.90"
Set-TextValue "E21" "  +4.07%  "
Set-TextValue "E22" "  +1.45%  "
Set-TextValue "E23" "  +3.07%  "
Set-TextValue "B24" "Fetch.AI"
Set-TextValue "C24" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D24" "2.40"
Set-TextValue "E24" "  +7.33%  "
Set-TextValue "B25" "InternetComputer(DFINITY)"
Set-TextValue "C25" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D25" "13.33"
Set-TextValue "E25" "  +3.34%  "
Set-TextValue "D26" "81.90"
Set-TextValue "E26" "  +0.97%  "
Set-TextValue "E27" "  -0.34%  "
Set-TextValue "D28" "9.97"
Set-TextValue "E28" "  +7.09%  "
Set-TextValue "D29" "7.49"
Set-TextValue "E29" "  +4.91%  "
Set-TextValue "D30" "2.72"
Set-TextValue "E30" "  +1.46%  "
Set-TextValue "B31" "ImmutableX"
Set-TextValue "C31" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D31" "2.22"
Set-TextValue "E31" "  +0.64%  "
Set-TextValue "B32" "FirstDigitalUSD"
Set-TextValue "C32" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D32" "1.00"
Set-TextValue "E32" "  -0.04%  "
Set-TextValue "E33" "  +7.06%  "
Set-TextValue "D34" "27.66"
Set-TextValue "E34" "  +2.64%  "
Set-TextValue "D35" "0.0₃0860"
Set-TextValue "E35" "  +0.68%  "
Set-TextValue "E36" "  +1.78%  "
Set-TextValue "D37" "3.41"
Set-TextValue "E37" "  +2.94%  "
Set-TextValue "D38" "6.19"
Set-TextValue "E38" "  +3.09%  "
Set-TextValue "D39" "2.28"
Set-TextValue "E39" "  -1.43%  "
Set-TextValue "E40" "  +6.59%  "
Set-TextValue "D41" "460.08"
Set-TextValue "E41" "  +4.79%  "
Set-TextValue "E42" "  -0.80%  "
Set-TextValue "D43" "0.291"
Set-TextValue "E43" "  +1.17%  "
Set-TextValue "D44" "0.0373"
Set-TextValue "E44" "  +0.48%  "
Set-TextValue "D45" "2.861.94"
Set-TextValue "E45" "  -1.72%  "
Set-TextValue "E46" "  +2.79%  "
Set-TextValue "D47" "39.45"
Set-TextValue "E47" "  -0.34%  "
Set-TextValue "D48" "130.89"
Set-TextValue "E48" "  +3.46%  "
Set-TextValue "D49" "25.62"
Set-TextValue "E49" "  +6.47%  "
Set-TextValue "D50" "2.30"
Set-TextValue "E50" "  +4.99%  "
